$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Sending cluster) changes for all data rows: Resolving-Mac -> ECs
$ws.Range("A2:A5").Value = "ECs"

# Row 2 (Target cluster = ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03556333333333333
$ws.Range("H2").Value = 0.10669
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 0.1349696900144444
$ws.Range("R2").Value = 1.21472721013
$ws.Range("S2").Value = 0.01044213755712683
$ws.Range("T2").Value = 0.01044213755712683

# Row 3 (Target cluster = FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03556333333333333
$ws.Range("H3").Value = 0.10669
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 8.655275991248889
$ws.Range("R3").Value = 77.89748392124001
$ws.Range("S3").Value = 0.6696287328350964
$ws.Range("T3").Value = 0.6696287328350964

# Row 4 (Target cluster = MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03556333333333333
$ws.Range("H4").Value = 0.10669
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 1.059836090663333
$ws.Range("R4").Value = 9.538524815969998
$ws.Range("S4").Value = 0.08199584844219236
$ws.Range("T4").Value = 0.08199584844219235

# Row 5 (Target cluster = Resolving-Mac)
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03556333333333333
$ws.Range("H5").Value = 0.10669
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 3.075403003202222
$ws.Range("R5").Value = 27.67862702882
$ws.Range("S5").Value = 0.2379332811655844
$ws.Range("T5").Value = 0.2379332811655844
